# Generate Report for Handoff
# The handoff/handback run picked up a new source file (new GUID) and new
# handoff timestamps; update the report cells + hyperlink display text to match.

$wb = $excel.ActiveWorkbook

$oldGuid = "22d7fe63-54eb-4fcb-88b0-bdab1e155ddf"
$newGuid = "868e5120-f9cc-45ae-a154-11922f54e8ed"

$linkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc2d46d5f0952b91d7005d47cd7940ffb41a5733/e2e/$oldGuid.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkTarget, [System.Type]::Missing, [System.Type]::Missing, "e2e\$newGuid.md") | Out-Null
$wsOverview.Range("G2").Value = "2016-08-24 19:07:34"

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $linkTarget, [System.Type]::Missing, [System.Type]::Missing, "$newGuid.md") | Out-Null
$wsZhCn.Range("G2").Value = "$newGuid.1349cc6be2b0898fd11a9dfe617f7cc85d069806.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-24 19:07:29"

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $linkTarget, [System.Type]::Missing, [System.Type]::Missing, "$newGuid.md") | Out-Null
$wsDeDe.Range("G2").Value = "$newGuid.1349cc6be2b0898fd11a9dfe617f7cc85d069806.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-24 19:07:34"
